$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2000
$ws.Range("J32").Value = 2000
$ws.Range("L32").Value = 2000
$ws.Range("N32").Value = -2652
$ws.Range("H74").Value = 5897.1055
$ws.Range("I74").Value = 5613.6113
$ws.Range("K74").Value = 5613.6113
$ws.Range("M74").Value = -4677.6113
$ws.Range("H77").Value = 5897.1055
$ws.Range("I77").Value = 5613.6113
$ws.Range("K77").Value = 28068.0565
$ws.Range("M77").Value = -23388.0565
$ws.Range("H111").Value = 2408.1904
$ws.Range("I111").Value = 2523.375
$ws.Range("K111").Value = 7570.125
$ws.Range("M111").Value = -4503.125
$ws.Range("H131").Value = 4860
$ws.Range("I131").Value = 1898.375
$ws.Range("K131").Value = 5695.125
$ws.Range("M131").Value = -655.125
$ws.Range("H132").Value = 41671028
$ws.Range("I132").Value = 43482784
$ws.Range("K132").Value = 130448352
$ws.Range("M132").Value = -130445822
$ws.Range("H137").Value = 11112744
$ws.Range("I137").Value = 41667656
$ws.Range("J137").Value = 1866.909
$ws.Range("K137").Value = 125002968
$ws.Range("L137").Value = 5600.727000000001
$ws.Range("M137").Value = -125000418
$ws.Range("N137").Value = -10700.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 6999.5
$ws.Range("I19").Value = 6999.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 6999.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -6770.5
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 20744.9
$ws.Range("I32").Value = 20306.375
$ws.Range("K32").Value = 20306.375
$ws.Range("M32").Value = -20019.375
$ws.Range("H45").Value = 1395.6
$ws.Range("I45").Value = 1395.6
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1395.6
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1018.6
$ws.Range("N45").ClearContents()
$ws.Range("H74").Value = 2142.2188
$ws.Range("I74").Value = 1253.7273
$ws.Range("J74").Value = 4096.9
$ws.Range("K74").Value = 1253.7273
$ws.Range("L74").Value = 4096.9
$ws.Range("M74").Value = -379.7273
$ws.Range("N74").Value = -5844.9
$ws.Range("H77").Value = 2142.2188
$ws.Range("I77").Value = 1253.7273
$ws.Range("J77").Value = 4096.9
$ws.Range("K77").Value = 6268.636500000001
$ws.Range("L77").Value = 20484.5
$ws.Range("M77").Value = -1900.636500000001
$ws.Range("N77").Value = -29220.5
$ws.Range("H132").Value = 3795.6843
$ws.Range("I132").Value = 2988.9033
$ws.Range("K132").Value = 8966.7099
$ws.Range("M132").Value = -6436.7099

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 16839.7
$ws.Range("I82").Value = 16839.7
$ws.Range("K82").Value = 16839.7
$ws.Range("M82").Value = -16456.7
$ws.Range("H85").Value = 16839.7
$ws.Range("I85").Value = 16839.7
$ws.Range("K85").Value = 16839.7
$ws.Range("M85").Value = -15513.7
$ws.Range("H94").Value = 863.9375
$ws.Range("I94").Value = 862.6923
$ws.Range("K94").Value = 862.6923
$ws.Range("M94").Value = -411.6923
$ws.Range("H107").Value = 2164.818
$ws.Range("I107").Value = 2080
$ws.Range("K107").Value = 2080
$ws.Range("M107").Value = -160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9453.786
$ws.Range("I86").Value = 10079.4
$ws.Range("J86").Value = 7889.75
$ws.Range("K86").Value = 10079.4
$ws.Range("L86").Value = 7889.75
$ws.Range("M86").Value = -8956.4
$ws.Range("N86").Value = -10135.75
$ws.Range("H89").Value = 9453.786
$ws.Range("I89").Value = 10079.4
$ws.Range("J89").Value = 7889.75
$ws.Range("K89").Value = 50397
$ws.Range("L89").Value = 39448.75
$ws.Range("M89").Value = -44781
$ws.Range("N89").Value = -50680.75
$ws.Range("H93").Value = 22703.5
$ws.Range("I93").Value = 22703.5
$ws.Range("K93").Value = 22703.5
$ws.Range("M93").Value = -20831.5
$ws.Range("H99").Value = 3072.8462
$ws.Range("I99").Value = 2730.6667
$ws.Range("J99").Value = 3842.75
$ws.Range("K99").Value = 2730.6667
$ws.Range("L99").Value = 3842.75
$ws.Range("M99").Value = -1232.6667
$ws.Range("N99").Value = -6838.75
$ws.Range("H126").Value = 3072.8462
$ws.Range("I126").Value = 2730.6667
$ws.Range("J126").Value = 3842.75
$ws.Range("K126").Value = 8192.000100000001
$ws.Range("L126").Value = 11528.25
$ws.Range("M126").Value = -5722.000100000001
$ws.Range("N126").Value = -16468.25
$ws.Range("H134").Value = 2030.421
$ws.Range("I134").Value = 2183.4666
$ws.Range("J134").Value = 1456.5
$ws.Range("K134").Value = 6550.399800000001
$ws.Range("L134").Value = 4369.5
$ws.Range("M134").Value = -4015.399800000001
$ws.Range("N134").Value = -9439.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3438164.5
$ws.Range("I4").Value = 1875703.8
$ws.Range("K4").Value = 5627111.4
$ws.Range("M4").Value = -5626999.4
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 9000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -9566

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5332
$ws.Range("H126").Value = 3960.5862
$ws.Range("I126").Value = 4062.15
$ws.Range("J126").Value = 3734.889
$ws.Range("K126").Value = 12186.45
$ws.Range("L126").Value = 11204.667
$ws.Range("M126").Value = -9716.450000000001
$ws.Range("N126").Value = -16144.667
$ws.Range("H132").Value = 9011674
$ws.Range("I132").Value = 2477.1072
$ws.Range("K132").Value = 7431.321599999999
$ws.Range("M132").Value = -4901.321599999999
$ws.Range("H136").Value = 63159
$ws.Range("J136").Value = 63159
$ws.Range("L136").Value = 189477
$ws.Range("N136").Value = -194577

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H122").Value = 3080
$ws.Range("I122").Value = 3080
$ws.Range("K122").Value = 9240
$ws.Range("M122").Value = -6790
$ws.Range("H132").Value = 2434.739
$ws.Range("J132").Value = 2449.9167
$ws.Range("L132").Value = 7349.750100000001
$ws.Range("N132").Value = -12409.7501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2507.4285
$ws.Range("I96").Value = 2507.4285
$ws.Range("K96").Value = 2507.4285
$ws.Range("M96").Value = -1134.4285
$ws.Range("H122").Value = 1887.3889
$ws.Range("I122").Value = 1887.3889
$ws.Range("K122").Value = 5662.1667
$ws.Range("M122").Value = -3212.1667
$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140
$ws.Range("H136").Value = 5963.0435
$ws.Range("I136").Value = 7745.0586
$ws.Range("K136").Value = 23235.1758
$ws.Range("M136").Value = -20685.1758
